# Remove the "Ministério da Saúde - Saúde A-Z" block (rows 12-18) from
# the worksheet. All rows below shift up by 7, matching the target
# layout where the sheet ends at row 28 instead of row 35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("12:18").Delete()
